$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 4110.1665
$ws.Cells.Item(86, 9).Value = 2705.5
$ws.Cells.Item(86, 11).Value = 2705.5
$ws.Cells.Item(86, 13).Value = -1582.5

$ws.Cells.Item(89, 8).Value = 4110.1665
$ws.Cells.Item(89, 9).Value = 2705.5
$ws.Cells.Item(89, 11).Value = 13527.5
$ws.Cells.Item(89, 13).Value = -7911.5

$ws.Cells.Item(112, 8).Value = 1065.2222
$ws.Cells.Item(112, 10).Value = 1085.5
$ws.Cells.Item(112, 12).Value = 3256.5
$ws.Cells.Item(112, 14).Value = -5472.5

$ws.Cells.Item(116, 8).Value = 198444.33
$ws.Cells.Item(116, 9).Value = 291666.5
$ws.Cells.Item(116, 10).Value = 12000
$ws.Cells.Item(116, 11).Value = 291666.5
$ws.Cells.Item(116, 12).Value = 12000
$ws.Cells.Item(116, 13).Value = -288224.5
$ws.Cells.Item(116, 14).Value = -18884

$ws.Cells.Item(133, 8).Value = 87662
$ws.Cells.Item(133, 10).Value = 87662
$ws.Cells.Item(133, 12).Value = 87662
$ws.Cells.Item(133, 14).Value = -97782

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 940.07465
$ws.Cells.Item(32, 9).Value = 913.8280999999999
$ws.Cells.Item(32, 11).Value = 913.8280999999999
$ws.Cells.Item(32, 13).Value = -626.8280999999999

$ws.Cells.Item(132, 8).Value = 3159.4
$ws.Cells.Item(132, 9).Value = 2624.4866
$ws.Cells.Item(132, 11).Value = 7873.459800000001
$ws.Cells.Item(132, 13).Value = -5343.459800000001

$ws.Cells.Item(139, 8).Value = 87142.60000000001
$ws.Cells.Item(139, 10).Value = 87142.60000000001
$ws.Cells.Item(139, 12).Value = 87142.60000000001
$ws.Cells.Item(139, 14).Value = -97422.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2883.1667
$ws.Cells.Item(86, 9).Value = 2859.8
$ws.Cells.Item(86, 10).Value = 3000
$ws.Cells.Item(86, 11).Value = 2859.8
$ws.Cells.Item(86, 12).Value = 3000
$ws.Cells.Item(86, 13).Value = -1736.8
$ws.Cells.Item(86, 14).Value = -5246

$ws.Cells.Item(89, 8).Value = 2883.1667
$ws.Cells.Item(89, 9).Value = 2859.8
$ws.Cells.Item(89, 10).Value = 3000
$ws.Cells.Item(89, 11).Value = 14299
$ws.Cells.Item(89, 12).Value = 15000
$ws.Cells.Item(89, 13).Value = -8683
$ws.Cells.Item(89, 14).Value = -26232

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 509.58823
$ws.Cells.Item(107, 10).Value = 923.8333
$ws.Cells.Item(107, 12).Value = 923.8333
$ws.Cells.Item(107, 14).Value = -4763.8333

$ws.Cells.Item(108, 8).Value = 63788.5

$ws.Cells.Item(111, 8).Value = 0
$ws.Cells.Item(111, 10).Value = 0
$ws.Cells.Item(111, 12).Value = 0
$ws.Cells.Item(111, 14).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 12383
$ws.Cells.Item(3, 9).Value = 2074.5
$ws.Cells.Item(3, 10).Value = 33000
$ws.Cells.Item(3, 11).Value = 6223.5
$ws.Cells.Item(3, 12).Value = 99000
$ws.Cells.Item(3, 13).Value = -6111.5
$ws.Cells.Item(3, 14).Value = -99224

$ws.Cells.Item(4, 8).Value = 32037484
$ws.Cells.Item(4, 9).Value = 38963404
$ws.Cells.Item(4, 11).Value = 116890212
$ws.Cells.Item(4, 13).Value = -116890100

$ws.Cells.Item(29, 8).Value = 1681.5
$ws.Cells.Item(29, 9).Value = 300
$ws.Cells.Item(29, 10).Value = 2372.25
$ws.Cells.Item(29, 11).Value = 900
$ws.Cells.Item(29, 12).Value = 7116.75
$ws.Cells.Item(29, 13).Value = -623
$ws.Cells.Item(29, 14).Value = -7670.75

$ws.Cells.Item(34, 8).Value = 935.2353000000001
$ws.Cells.Item(34, 10).Value = 9000
$ws.Cells.Item(34, 12).Value = 27000
$ws.Cells.Item(34, 14).Value = -27168

$ws.Cells.Item(37, 8).Value = 97124.875
$ws.Cells.Item(37, 10).Value = 97124.875
$ws.Cells.Item(37, 12).Value = 291374.625
$ws.Cells.Item(37, 14).Value = -291598.625

$ws.Cells.Item(46, 8).Value = 217.66667
$ws.Cells.Item(46, 9).Value = 217.66667
$ws.Cells.Item(46, 11).Value = 653.00001
$ws.Cells.Item(46, 13).Value = -562.00001

$ws.Cells.Item(141, 8).Value = 2681.9167
$ws.Cells.Item(141, 9).Value = 2598.4546
$ws.Cells.Item(141, 10).Value = 3600
$ws.Cells.Item(141, 11).Value = 7795.3638
$ws.Cells.Item(141, 12).Value = 10800
$ws.Cells.Item(141, 13).Value = -2615.3638
$ws.Cells.Item(141, 14).Value = -21160

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 14).ClearContents()

$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 13).ClearContents()

$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 13).ClearContents()

$ws.Cells.Item(97, 8).Value = 2283.1333
$ws.Cells.Item(97, 9).Value = 2363.6
$ws.Cells.Item(97, 10).Value = 2122.2
$ws.Cells.Item(97, 11).Value = 2363.6
$ws.Cells.Item(97, 12).Value = 2122.2
$ws.Cells.Item(97, 13).Value = -1867.6
$ws.Cells.Item(97, 14).Value = -3114.2

$ws.Cells.Item(102, 8).Value = 3453.6667
$ws.Cells.Item(102, 9).Value = 3260.375
$ws.Cells.Item(102, 11).Value = 3260.375
$ws.Cells.Item(102, 13).Value = -1638.375

$ws.Cells.Item(132, 8).Value = 24403612
$ws.Cells.Item(132, 9).Value = 35727156
$ws.Cells.Item(132, 11).Value = 107181468
$ws.Cells.Item(132, 13).Value = -107178938

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 3418.6843
$ws.Cells.Item(68, 9).Value = 2380.3076
$ws.Cells.Item(68, 11).Value = 2380.3076
$ws.Cells.Item(68, 13).Value = -1631.3076

$ws.Cells.Item(71, 8).Value = 3418.6843
$ws.Cells.Item(71, 9).Value = 2380.3076
$ws.Cells.Item(71, 11).Value = 11901.538
$ws.Cells.Item(71, 13).Value = -8157.538

$ws.Cells.Item(116, 8).Value = 0
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 14).ClearContents()

$ws.Cells.Item(120, 8).Value = 97581.664
$ws.Cells.Item(120, 10).Value = 97581.664
$ws.Cells.Item(120, 12).Value = 97581.664
$ws.Cells.Item(120, 14).Value = -107257.664

$ws.Cells.Item(122, 8).Value = 5056.095
$ws.Cells.Item(122, 9).Value = 4598.706
$ws.Cells.Item(122, 11).Value = 13796.118
$ws.Cells.Item(122, 13).Value = -11346.118

$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 12).Value = 0
$ws.Cells.Item(133, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4865
$ws.Cells.Item(62, 9).Value = 5097.5
$ws.Cells.Item(62, 10).Value = 4400
$ws.Cells.Item(62, 11).Value = 5097.5
$ws.Cells.Item(62, 12).Value = 4400
$ws.Cells.Item(62, 13).Value = -4473.5
$ws.Cells.Item(62, 14).Value = -5648

$ws.Cells.Item(65, 8).Value = 4865
$ws.Cells.Item(65, 9).Value = 5097.5
$ws.Cells.Item(65, 10).Value = 4400
$ws.Cells.Item(65, 11).Value = 25487.5
$ws.Cells.Item(65, 12).Value = 22000
$ws.Cells.Item(65, 13).Value = -22367.5
$ws.Cells.Item(65, 14).Value = -28240

$ws.Cells.Item(81, 8).Value = 5829.273
$ws.Cells.Item(81, 9).Value = 6282.2
$ws.Cells.Item(81, 11).Value = 12564.4
$ws.Cells.Item(81, 13).Value = -11503.4

$ws.Cells.Item(84, 8).Value = 5829.273
$ws.Cells.Item(84, 9).Value = 6282.2
$ws.Cells.Item(84, 11).Value = 62822
$ws.Cells.Item(84, 13).Value = -57518

$ws.Cells.Item(107, 8).Value = 4364.5713
$ws.Cells.Item(107, 9).Value = 2489.4285
$ws.Cells.Item(107, 10).Value = 6239.7144
$ws.Cells.Item(107, 11).Value = 7468.2855
$ws.Cells.Item(107, 12).Value = 18719.1432
$ws.Cells.Item(107, 13).Value = -5548.2855
$ws.Cells.Item(107, 14).Value = -22559.1432
